$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the status value for the "Front Right" row (B2): was "B", should be "O"
$ws.Range("B2").Value = "O"

# Update the active selection to B2, matching the saved view state
$ws.Range("B2").Select()
